# Update "想去人数" (F column) and "最低票价" (G column) figures on the
# "展览" and "全部类型" worksheets, which hold identical data.
#
# Changes (row -> old -> new):
#   F2:  4911 -> 4912
#   F5:  801  -> 802
#   G5:  52.1 -> 55
#   F7:  1279 -> 1282
#   F8:  136  -> 137
#   F10: 214  -> 215
#   F13: 166  -> 167
#   F15: 4309 -> 4319
#   F16: 6607 -> 6610
#   F20: 555  -> 556
#   F23: 427  -> 428
#   F24: 59   -> 61
#   F25: 36   -> 38
#   F26: 2650 -> 2656
#   F27: 571  -> 572
#   F32: 389  -> 390
#   F33: 206  -> 207
#   F34: 24   -> 25
#   F35: 1598 -> 1601
#   F38: 112  -> 114
#   F40: 519  -> 521
#   F43: 82   -> 83
#   F44: 610  -> 611

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 4912
    "F5"  = 802
    "G5"  = 55
    "F7"  = 1282
    "F8"  = 137
    "F10" = 215
    "F13" = 167
    "F15" = 4319
    "F16" = 6610
    "F20" = 556
    "F23" = 428
    "F24" = 61
    "F25" = 38
    "F26" = 2656
    "F27" = 572
    "F32" = 390
    "F33" = 207
    "F34" = 25
    "F35" = 1601
    "F38" = 114
    "F40" = 521
    "F43" = 83
    "F44" = 611
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
